# Daily COVID data refresh for the "Pais" ranking sheet.
# Row 1 holds the "last updated" timestamp string; row 3 is the (unchanged)
# header row; rows 4.. hold one country per row, pre-sorted descending by
# total cases (col B). New counts shift a few countries past their former
# neighbour, so rows are rewritten as a whole block to keep the table in order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 18 de Septiembre de 2020 a las 01:14'

$rows = 216
$arr = New-Object "object[,]" $rows,8
$arr[0,0] = 'Estados Unidos'; $arr[0,1] = 6869238; $arr[0,2] = 40937; $arr[0,3] = 4145890; $arr[0,4] = 2521225; $arr[0,5] = 0; $arr[0,6] = 789; $arr[0,7] = 202123
$arr[1,0] = 'India'; $arr[1,1] = 5212686; $arr[1,2] = 96793; $arr[1,3] = 4109828; $arr[1,4] = 1018454; $arr[1,5] = 0; $arr[1,6] = 1174; $arr[1,7] = 84404
$arr[2,0] = 'Brasil'; $arr[2,1] = 4457443; $arr[2,2] = 35757; $arr[2,3] = 3753082; $arr[2,4] = 569330; $arr[2,5] = 0; $arr[2,6] = 857; $arr[2,7] = 135031
$arr[3,0] = 'Rusia'; $arr[3,1] = 1085281; $arr[3,2] = 5762; $arr[3,3] = 895868; $arr[3,4] = 170352; $arr[3,5] = 0; $arr[3,6] = 144; $arr[3,7] = 19061
$arr[4,0] = 'Peru'; $arr[4,1] = 744400; $arr[4,2] = 0; $arr[4,3] = 587717; $arr[4,4] = 125632; $arr[4,5] = 0; $arr[4,6] = 0; $arr[4,7] = 31051
$arr[5,0] = 'Colombia'; $arr[5,1] = 743945; $arr[5,2] = 7568; $arr[5,3] = 615457; $arr[5,4] = 104823; $arr[5,5] = 0; $arr[5,6] = 187; $arr[5,7] = 23665
$arr[6,0] = 'Mexico'; $arr[6,1] = 680931; $arr[6,2] = 4444; $arr[6,3] = 485024; $arr[6,4] = 123929; $arr[6,5] = 0; $arr[6,6] = 300; $arr[6,7] = 71978
$arr[7,0] = 'Sudafrica'; $arr[7,1] = 655572; $arr[7,2] = 2128; $arr[7,3] = 585303; $arr[7,4] = 54497; $arr[7,5] = 0; $arr[7,6] = 67; $arr[7,7] = 15772
$arr[8,0] = 'España'; $arr[8,1] = 625651; $arr[8,2] = 11291; $arr[8,3] = 0; $arr[8,4] = 0; $arr[8,5] = 0; $arr[8,6] = 162; $arr[8,7] = 30405
$arr[9,0] = 'Argentina'; $arr[9,1] = 601713; $arr[9,2] = 12701; $arr[9,3] = 456347; $arr[9,4] = 132906; $arr[9,5] = 0; $arr[9,6] = 344; $arr[9,7] = 12460
$arr[10,0] = 'Chile'; $arr[10,1] = 441150; $arr[10,2] = 1863; $arr[10,3] = 413928; $arr[10,4] = 15080; $arr[10,5] = 0; $arr[10,6] = 84; $arr[10,7] = 12142
$arr[11,0] = 'Francia'; $arr[11,1] = 415481; $arr[11,2] = 10593; $arr[11,3] = 90840; $arr[11,4] = 293546; $arr[11,5] = 0; $arr[11,6] = 50; $arr[11,7] = 31095
$arr[12,0] = 'Iran'; $arr[12,1] = 413149; $arr[12,2] = 2815; $arr[12,3] = 353848; $arr[12,4] = 35493; $arr[12,5] = 0; $arr[12,6] = 176; $arr[12,7] = 23808
$arr[13,0] = 'Reino Unido'; $arr[13,1] = 381614; $arr[13,2] = 3395; $arr[13,3] = 0; $arr[13,4] = 0; $arr[13,5] = 0; $arr[13,6] = 21; $arr[13,7] = 41705
$arr[14,0] = 'Banglades'; $arr[14,1] = 344264; $arr[14,2] = 1593; $arr[14,3] = 250412; $arr[14,4] = 88993; $arr[14,5] = 0; $arr[14,6] = 36; $arr[14,7] = 4859
$arr[15,0] = 'Arabia Saudita'; $arr[15,1] = 328144; $arr[15,2] = 593; $arr[15,3] = 307207; $arr[15,4] = 16538; $arr[15,5] = 0; $arr[15,6] = 30; $arr[15,7] = 4399
$arr[16,0] = 'Irak'; $arr[16,1] = 307385; $arr[16,2] = 4326; $arr[16,3] = 241100; $arr[16,4] = 57953; $arr[16,5] = 0; $arr[16,6] = 84; $arr[16,7] = 8332
$arr[17,0] = 'Pakistan'; $arr[17,1] = 303634; $arr[17,2] = 545; $arr[17,3] = 291169; $arr[17,4] = 6066; $arr[17,5] = 0; $arr[17,6] = 6; $arr[17,7] = 6399
$arr[18,0] = 'Turquia'; $arr[18,1] = 298039; $arr[18,2] = 1648; $arr[18,3] = 263745; $arr[18,4] = 26979; $arr[18,5] = 0; $arr[18,6] = 66; $arr[18,7] = 7315
$arr[19,0] = 'Italia'; $arr[19,1] = 293025; $arr[19,2] = 1585; $arr[19,3] = 215954; $arr[19,4] = 41413; $arr[19,5] = 0; $arr[19,6] = 13; $arr[19,7] = 35658
$arr[20,0] = 'Filipinas'; $arr[20,1] = 276289; $arr[20,2] = 3375; $arr[20,3] = 208096; $arr[20,4] = 63408; $arr[20,5] = 0; $arr[20,6] = 53; $arr[20,7] = 4785
$arr[21,0] = 'Alemania'; $arr[21,1] = 269035; $arr[21,2] = 2170; $arr[21,3] = 239100; $arr[21,4] = 20480; $arr[21,5] = 0; $arr[21,6] = 6; $arr[21,7] = 9455
$arr[22,0] = 'Indonesia'; $arr[22,1] = 232628; $arr[22,2] = 3635; $arr[22,3] = 166686; $arr[22,4] = 56720; $arr[22,5] = 0; $arr[22,6] = 122; $arr[22,7] = 9222
$arr[23,0] = 'Israel'; $arr[23,1] = 175256; $arr[23,2] = 4791; $arr[23,3] = 126329; $arr[23,4] = 47758; $arr[23,5] = 0; $arr[23,6] = 8; $arr[23,7] = 1169
$arr[24,0] = 'Ucrania'; $arr[24,1] = 166244; $arr[24,2] = 3584; $arr[24,3] = 73913; $arr[24,4] = 88931; $arr[24,5] = 0; $arr[24,6] = 60; $arr[24,7] = 3400
$arr[25,0] = 'Canada'; $arr[25,1] = 140539; $arr[25,2] = 792; $arr[25,3] = 122836; $arr[25,4] = 8504; $arr[25,5] = 0; $arr[25,6] = 6; $arr[25,7] = 9199
$arr[26,0] = 'Bolivia'; $arr[26,1] = 128872; $arr[26,2] = 586; $arr[26,3] = 87031; $arr[26,4] = 34363; $arr[26,5] = 0; $arr[26,6] = 31; $arr[26,7] = 7478
$arr[27,0] = 'Catar'; $arr[27,1] = 122693; $arr[27,2] = 244; $arr[27,3] = 119613; $arr[27,4] = 2872; $arr[27,5] = 0; $arr[27,6] = 0; $arr[27,7] = 208
$arr[28,0] = 'Ecuador'; $arr[28,1] = 122257; $arr[28,2] = 732; $arr[28,3] = 97063; $arr[28,4] = 14165; $arr[28,5] = 0; $arr[28,6] = 33; $arr[28,7] = 11029
$arr[29,0] = 'Rumania'; $arr[29,1] = 108690; $arr[29,2] = 1679; $arr[29,3] = 43244; $arr[29,4] = 61134; $arr[29,5] = 0; $arr[29,6] = 27; $arr[29,7] = 4312
$arr[30,0] = 'Kazajistan'; $arr[30,1] = 107056; $arr[30,2] = 72; $arr[30,3] = 101455; $arr[30,4] = 3930; $arr[30,5] = 0; $arr[30,6] = 0; $arr[30,7] = 1671
$arr[31,0] = 'Republica Dominicana'; $arr[31,1] = 106136; $arr[31,2] = 615; $arr[31,3] = 79363; $arr[31,4] = 24751; $arr[31,5] = 0; $arr[31,6] = 13; $arr[31,7] = 2022
$arr[32,0] = 'Panama'; $arr[32,1] = 104138; $arr[32,2] = 672; $arr[32,3] = 77881; $arr[32,4] = 24044; $arr[32,5] = 0; $arr[32,6] = 15; $arr[32,7] = 2213
$arr[33,0] = 'Egipto'; $arr[33,1] = 101641; $arr[33,2] = 141; $arr[33,3] = 87158; $arr[33,4] = 8768; $arr[33,5] = 0; $arr[33,6] = 19; $arr[33,7] = 5715
$arr[34,0] = 'Kuwait'; $arr[34,1] = 97824; $arr[34,2] = 825; $arr[34,3] = 87911; $arr[34,4] = 9338; $arr[34,5] = 0; $arr[34,6] = 4; $arr[34,7] = 575
$arr[35,0] = 'Belgica'; $arr[35,1] = 95948; $arr[35,2] = 1153; $arr[35,3] = 18810; $arr[35,4] = 67203; $arr[35,5] = 0; $arr[35,6] = 5; $arr[35,7] = 9935
$arr[36,0] = 'Marruecos'; $arr[36,1] = 94504; $arr[36,2] = 2488; $arr[36,3] = 74930; $arr[36,4] = 17860; $arr[36,5] = 0; $arr[36,6] = 28; $arr[36,7] = 1714
$arr[37,0] = 'Oman'; $arr[37,1] = 91753; $arr[37,2] = 557; $arr[37,3] = 84648; $arr[37,4] = 6287; $arr[37,5] = 0; $arr[37,6] = 13; $arr[37,7] = 818
$arr[38,0] = 'Paises Bajos'; $arr[38,1] = 88073; $arr[38,2] = 1753; $arr[38,3] = 0; $arr[38,4] = 0; $arr[38,5] = 0; $arr[38,6] = 6; $arr[38,7] = 6266
$arr[39,0] = 'Suecia'; $arr[39,1] = 87885; $arr[39,2] = 0; $arr[39,3] = 0; $arr[39,4] = 0; $arr[39,5] = 0; $arr[39,6] = 5; $arr[39,7] = 5864
$arr[40,0] = 'China'; $arr[40,1] = 85223; $arr[40,2] = 9; $arr[40,3] = 80448; $arr[40,4] = 141; $arr[40,5] = 0; $arr[40,6] = 0; $arr[40,7] = 4634
$arr[41,0] = 'Guatemala'; $arr[41,1] = 83664; $arr[41,2] = 740; $arr[41,3] = 73260; $arr[41,4] = 7368; $arr[41,5] = 0; $arr[41,6] = 27; $arr[41,7] = 3036
$arr[42,0] = 'Emiratos Arabes Unidos'; $arr[42,1] = 82568; $arr[42,2] = 786; $arr[42,3] = 72117; $arr[42,4] = 10049; $arr[42,5] = 0; $arr[42,6] = 0; $arr[42,7] = 402
$arr[43,0] = 'Japon'; $arr[43,1] = 77009; $arr[43,2] = 561; $arr[43,3] = 69253; $arr[43,4] = 6283; $arr[43,5] = 0; $arr[43,6] = 12; $arr[43,7] = 1473
$arr[44,0] = 'Polonia'; $arr[44,1] = 76571; $arr[44,2] = 837; $arr[44,3] = 62725; $arr[44,4] = 11593; $arr[44,5] = 0; $arr[44,6] = 16; $arr[44,7] = 2253
$arr[45,0] = 'Bielorrusia'; $arr[45,1] = 74987; $arr[45,2] = 224; $arr[45,3] = 72967; $arr[45,4] = 1249; $arr[45,5] = 0; $arr[45,6] = 4; $arr[45,7] = 771
$arr[46,0] = 'Honduras'; $arr[46,1] = 69660; $arr[46,2] = 1040; $arr[46,3] = 19983; $arr[46,4] = 47575; $arr[46,5] = 0; $arr[46,6] = 15; $arr[46,7] = 2102
$arr[47,0] = 'Etiopia'; $arr[47,1] = 66913; $arr[47,2] = 689; $arr[47,3] = 27085; $arr[47,4] = 38768; $arr[47,5] = 0; $arr[47,6] = 15; $arr[47,7] = 1060
$arr[48,0] = 'Portugal'; $arr[48,1] = 66396; $arr[48,2] = 770; $arr[48,3] = 44794; $arr[48,4] = 19714; $arr[48,5] = 0; $arr[48,6] = 10; $arr[48,7] = 1888
$arr[49,0] = 'Venezuela'; $arr[49,1] = 63416; $arr[49,2] = 0; $arr[49,3] = 51274; $arr[49,4] = 11631; $arr[49,5] = 0; $arr[49,6] = 0; $arr[49,7] = 511
$arr[50,0] = 'Barein'; $arr[50,1] = 63189; $arr[50,2] = 705; $arr[50,3] = 56087; $arr[50,4] = 6885; $arr[50,5] = 0; $arr[50,6] = 1; $arr[50,7] = 217
$arr[51,0] = 'Costa Rica'; $arr[51,1] = 60818; $arr[51,2] = 1302; $arr[51,3] = 22662; $arr[51,4] = 37490; $arr[51,5] = 0; $arr[51,6] = 17; $arr[51,7] = 666
$arr[52,0] = 'Nepal'; $arr[52,1] = 59573; $arr[52,2] = 1246; $arr[52,3] = 42949; $arr[52,4] = 16241; $arr[52,5] = 0; $arr[52,6] = 4; $arr[52,7] = 383
$arr[53,0] = 'Singapur'; $arr[53,1] = 57532; $arr[53,2] = 18; $arr[53,3] = 57039; $arr[53,4] = 466; $arr[53,5] = 0; $arr[53,6] = 0; $arr[53,7] = 27
$arr[54,0] = 'Nigeria'; $arr[54,1] = 56735; $arr[54,2] = 131; $arr[54,3] = 48092; $arr[54,4] = 7550; $arr[54,5] = 0; $arr[54,6] = 2; $arr[54,7] = 1093
$arr[55,0] = 'Uzbekistan'; $arr[55,1] = 49627; $arr[55,2] = 612; $arr[55,3] = 45970; $arr[55,4] = 3244; $arr[55,5] = 0; $arr[55,6] = 6; $arr[55,7] = 413
$arr[56,0] = 'Argelia'; $arr[56,1] = 49194; $arr[56,2] = 228; $arr[56,3] = 34675; $arr[56,4] = 12865; $arr[56,5] = 0; $arr[56,6] = 9; $arr[56,7] = 1654
$arr[57,0] = 'Suiza'; $arr[57,1] = 48795; $arr[57,2] = 530; $arr[57,3] = 39900; $arr[57,4] = 6853; $arr[57,5] = 0; $arr[57,6] = 3; $arr[57,7] = 2042
$arr[58,0] = 'Armenia'; $arr[58,1] = 46671; $arr[58,2] = 295; $arr[58,3] = 42231; $arr[58,4] = 3515; $arr[58,5] = 0; $arr[58,6] = 2; $arr[58,7] = 925
$arr[59,0] = 'Ghana'; $arr[59,1] = 45714; $arr[59,2] = 59; $arr[59,3] = 44896; $arr[59,4] = 524; $arr[59,5] = 0; $arr[59,6] = 0; $arr[59,7] = 294
$arr[60,0] = 'Kirguistan'; $arr[60,1] = 45153; $arr[60,2] = 81; $arr[60,3] = 41317; $arr[60,4] = 2773; $arr[60,5] = 0; $arr[60,6] = 0; $arr[60,7] = 1063
$arr[61,0] = 'Moldavia'; $arr[61,1] = 44983; $arr[61,2] = 622; $arr[61,3] = 33239; $arr[61,4] = 10574; $arr[61,5] = 0; $arr[61,6] = 11; $arr[61,7] = 1170
$arr[62,0] = 'Chequia'; $arr[62,1] = 42739; $arr[62,2] = 1707; $arr[62,3] = 23321; $arr[62,4] = 18930; $arr[62,5] = 0; $arr[62,6] = 6; $arr[62,7] = 488
$arr[63,0] = 'Afganistan'; $arr[63,1] = 38872; $arr[63,2] = 17; $arr[63,3] = 32505; $arr[63,4] = 4931; $arr[63,5] = 0; $arr[63,6] = 0; $arr[63,7] = 1436
$arr[64,0] = 'Azerbaiyan'; $arr[64,1] = 38777; $arr[64,2] = 119; $arr[64,3] = 36289; $arr[64,4] = 1917; $arr[64,5] = 0; $arr[64,6] = 2; $arr[64,7] = 571
$arr[65,0] = 'Kenia'; $arr[65,1] = 36576; $arr[65,2] = 183; $arr[65,3] = 23611; $arr[65,4] = 12323; $arr[65,5] = 0; $arr[65,6] = 5; $arr[65,7] = 642
$arr[66,0] = 'Austria'; $arr[66,1] = 35853; $arr[66,2] = 780; $arr[66,3] = 28044; $arr[66,4] = 7051; $arr[66,5] = 0; $arr[66,6] = 0; $arr[66,7] = 758
$arr[67,0] = 'Estado de Palestina'; $arr[67,1] = 33843; $arr[67,2] = 837; $arr[67,3] = 23060; $arr[67,4] = 10539; $arr[67,5] = 0; $arr[67,6] = 1; $arr[67,7] = 244
$arr[68,0] = 'Serbia'; $arr[68,1] = 32695; $arr[68,2] = 82; $arr[68,3] = 31512; $arr[68,4] = 445; $arr[68,5] = 0; $arr[68,6] = 2; $arr[68,7] = 738
$arr[69,0] = 'Irlanda'; $arr[69,1] = 32023; $arr[69,2] = 224; $arr[69,3] = 23364; $arr[69,4] = 6870; $arr[69,5] = 0; $arr[69,6] = 1; $arr[69,7] = 1789
$arr[70,0] = 'Paraguay'; $arr[70,1] = 30419; $arr[70,2] = 0; $arr[70,3] = 15740; $arr[70,4] = 14113; $arr[70,5] = 0; $arr[70,6] = 0; $arr[70,7] = 566
$arr[71,0] = 'El Salvador'; $arr[71,1] = 27249; $arr[71,2] = 86; $arr[71,3] = 20392; $arr[71,4] = 6056; $arr[71,5] = 0; $arr[71,6] = 5; $arr[71,7] = 801
$arr[72,0] = 'Australia'; $arr[72,1] = 26813; $arr[72,2] = 34; $arr[72,3] = 23792; $arr[72,4] = 2189; $arr[72,5] = 0; $arr[72,6] = 8; $arr[72,7] = 832
$arr[73,0] = 'Libano'; $arr[73,1] = 26768; $arr[73,2] = 685; $arr[73,3] = 10217; $arr[73,4] = 16288; $arr[73,5] = 0; $arr[73,6] = 4; $arr[73,7] = 263
$arr[74,0] = 'Libia'; $arr[74,1] = 25822; $arr[74,2] = 886; $arr[74,3] = 13908; $arr[74,4] = 11505; $arr[74,5] = 0; $arr[74,6] = 15; $arr[74,7] = 409
$arr[75,0] = 'Bosnia y Herzegovina'; $arr[75,1] = 24605; $arr[75,2] = 394; $arr[75,3] = 17219; $arr[75,4] = 6639; $arr[75,5] = 0; $arr[75,6] = 11; $arr[75,7] = 747
$arr[76,0] = 'Corea del Sur'; $arr[76,1] = 22657; $arr[76,2] = 153; $arr[76,3] = 19543; $arr[76,4] = 2742; $arr[76,5] = 0; $arr[76,6] = 5; $arr[76,7] = 372
$arr[77,0] = 'Dinamarca'; $arr[77,1] = 21393; $arr[77,2] = 453; $arr[77,3] = 16918; $arr[77,4] = 3840; $arr[77,5] = 0; $arr[77,6] = 1; $arr[77,7] = 635
$arr[78,0] = 'Camerun'; $arr[78,1] = 20303; $arr[78,2] = 0; $arr[78,3] = 18837; $arr[78,4] = 1051; $arr[78,5] = 0; $arr[78,6] = 0; $arr[78,7] = 415
$arr[79,0] = 'Costa de Marfil'; $arr[79,1] = 19158; $arr[79,2] = 26; $arr[79,3] = 18330; $arr[79,4] = 708; $arr[79,5] = 0; $arr[79,6] = 0; $arr[79,7] = 120
$arr[80,0] = 'Bulgaria'; $arr[80,1] = 18544; $arr[80,2] = 154; $arr[80,3] = 13391; $arr[80,4] = 4404; $arr[80,5] = 0; $arr[80,6] = 10; $arr[80,7] = 749
$arr[81,0] = 'Republica de Macedonia'; $arr[81,1] = 16274; $arr[81,2] = 186; $arr[81,3] = 13635; $arr[81,4] = 1964; $arr[81,5] = 0; $arr[81,6] = 7; $arr[81,7] = 675
$arr[82,0] = 'Madagascar'; $arr[82,1] = 15925; $arr[82,2] = 54; $arr[82,3] = 14547; $arr[82,4] = 1162; $arr[82,5] = 0; $arr[82,6] = 1; $arr[82,7] = 216
$arr[83,0] = 'Hungria'; $arr[83,1] = 15170; $arr[83,2] = 710; $arr[83,3] = 4227; $arr[83,4] = 10280; $arr[83,5] = 0; $arr[83,6] = 9; $arr[83,7] = 663
$arr[84,0] = 'Senegal'; $arr[84,1] = 14618; $arr[84,2] = 50; $arr[84,3] = 10923; $arr[84,4] = 3395; $arr[84,5] = 0; $arr[84,6] = 1; $arr[84,7] = 300
$arr[85,0] = 'Grecia'; $arr[85,1] = 14400; $arr[85,2] = 359; $arr[85,3] = 3804; $arr[85,4] = 10271; $arr[85,5] = 0; $arr[85,6] = 9; $arr[85,7] = 325
$arr[86,0] = 'Croacia'; $arr[86,1] = 14279; $arr[86,2] = 250; $arr[86,3] = 11933; $arr[86,4] = 2108; $arr[86,5] = 0; $arr[86,6] = 2; $arr[86,7] = 238
$arr[87,0] = 'Zambia'; $arr[87,1] = 13928; $arr[87,2] = 41; $arr[87,3] = 13029; $arr[87,4] = 573; $arr[87,5] = 0; $arr[87,6] = 0; $arr[87,7] = 326
$arr[88,0] = 'Sudan'; $arr[88,1] = 13535; $arr[88,2] = 0; $arr[88,3] = 6759; $arr[88,4] = 5940; $arr[88,5] = 0; $arr[88,6] = 0; $arr[88,7] = 836
$arr[89,0] = 'Noruega'; $arr[89,1] = 12571; $arr[89,2] = 73; $arr[89,3] = 10371; $arr[89,4] = 1934; $arr[89,5] = 0; $arr[89,6] = 1; $arr[89,7] = 266
$arr[90,0] = 'Albania'; $arr[90,1] = 11948; $arr[90,2] = 132; $arr[90,3] = 6788; $arr[90,4] = 4813; $arr[90,5] = 0; $arr[90,6] = 4; $arr[90,7] = 347
$arr[91,0] = 'Consejo Danes para los Refugiados'; $arr[91,1] = 10442; $arr[91,2] = 28; $arr[91,3] = 9840; $arr[91,4] = 335; $arr[91,5] = 0; $arr[91,6] = 0; $arr[91,7] = 267
$arr[92,0] = 'Guinea'; $arr[92,1] = 10154; $arr[92,2] = 0; $arr[92,3] = 9612; $arr[92,4] = 479; $arr[92,5] = 0; $arr[92,6] = 0; $arr[92,7] = 63
$arr[93,0] = 'Namibia'; $arr[93,1] = 10078; $arr[93,2] = 114; $arr[93,3] = 7685; $arr[93,4] = 2285; $arr[93,5] = 0; $arr[93,6] = 0; $arr[93,7] = 108
$arr[94,0] = 'Malasia'; $arr[94,1] = 10052; $arr[94,2] = 21; $arr[94,3] = 9250; $arr[94,4] = 674; $arr[94,5] = 0; $arr[94,6] = 0; $arr[94,7] = 128
$arr[95,0] = 'Guayana Francesa'; $arr[95,1] = 9623; $arr[95,2] = 28; $arr[95,3] = 9267; $arr[95,4] = 291; $arr[95,5] = 0; $arr[95,6] = 0; $arr[95,7] = 65
$arr[96,0] = 'Maldivas'; $arr[96,1] = 9494; $arr[96,2] = 67; $arr[96,3] = 8033; $arr[96,4] = 1428; $arr[96,5] = 0; $arr[96,6] = 0; $arr[96,7] = 33
$arr[97,0] = 'Tayikistan'; $arr[97,1] = 9214; $arr[97,2] = 43; $arr[97,3] = 7988; $arr[97,4] = 1153; $arr[97,5] = 0; $arr[97,6] = 0; $arr[97,7] = 73
$arr[98,0] = 'Finlandia'; $arr[98,1] = 8799; $arr[98,2] = 49; $arr[98,3] = 7700; $arr[98,4] = 760; $arr[98,5] = 0; $arr[98,6] = 0; $arr[98,7] = 339
$arr[99,0] = 'Gabon'; $arr[99,1] = 8678; $arr[99,2] = 0; $arr[99,3] = 7827; $arr[99,4] = 798; $arr[99,5] = 0; $arr[99,6] = 0; $arr[99,7] = 53
$arr[100,0] = 'Tunez'; $arr[100,1] = 8570; $arr[100,2] = 470; $arr[100,3] = 2342; $arr[100,4] = 6095; $arr[100,5] = 0; $arr[100,6] = 4; $arr[100,7] = 133
$arr[101,0] = 'Haiti'; $arr[101,1] = 8556; $arr[101,2] = 15; $arr[101,3] = 6315; $arr[101,4] = 2021; $arr[101,5] = 0; $arr[101,6] = 0; $arr[101,7] = 220
$arr[102,0] = 'Zimbabue'; $arr[102,1] = 7633; $arr[102,2] = 35; $arr[102,3] = 5841; $arr[102,4] = 1568; $arr[102,5] = 0; $arr[102,6] = 0; $arr[102,7] = 224
$arr[103,0] = 'Luxemburgo'; $arr[103,1] = 7394; $arr[103,2] = 0; $arr[103,3] = 6593; $arr[103,4] = 677; $arr[103,5] = 0; $arr[103,6] = 0; $arr[103,7] = 124
$arr[104,0] = 'Mauritania'; $arr[104,1] = 7346; $arr[104,2] = 14; $arr[104,3] = 6865; $arr[104,4] = 320; $arr[104,5] = 0; $arr[104,6] = 0; $arr[104,7] = 161
$arr[105,0] = 'Montenegro'; $arr[105,1] = 7291; $arr[105,2] = 0; $arr[105,3] = 4764; $arr[105,4] = 2401; $arr[105,5] = 0; $arr[105,6] = 0; $arr[105,7] = 126
$arr[106,0] = 'Mozambique'; $arr[106,1] = 6161; $arr[106,2] = 167; $arr[106,3] = 3393; $arr[106,4] = 2729; $arr[106,5] = 0; $arr[106,6] = 0; $arr[106,7] = 39
$arr[107,0] = 'Eslovaquia'; $arr[107,1] = 6021; $arr[107,2] = 161; $arr[107,3] = 3288; $arr[107,4] = 2694; $arr[107,5] = 0; $arr[107,6] = 1; $arr[107,7] = 39
$arr[108,0] = 'Malaui'; $arr[108,1] = 5711; $arr[108,2] = 7; $arr[108,3] = 4000; $arr[108,4] = 1532; $arr[108,5] = 0; $arr[108,6] = 1; $arr[108,7] = 179
$arr[109,0] = 'Republica de Yibuti'; $arr[109,1] = 5399; $arr[109,2] = 0; $arr[109,3] = 5333; $arr[109,4] = 5; $arr[109,5] = 0; $arr[109,6] = 0; $arr[109,7] = 61
$arr[110,0] = 'Uganda'; $arr[110,1] = 5380; $arr[110,2] = 114; $arr[110,3] = 2489; $arr[110,4] = 2831; $arr[110,5] = 0; $arr[110,6] = 0; $arr[110,7] = 60
$arr[111,0] = 'Suazilandia'; $arr[111,1] = 5155; $arr[111,2] = 0; $arr[111,3] = 4418; $arr[111,4] = 636; $arr[111,5] = 0; $arr[111,6] = 0; $arr[111,7] = 101
$arr[112,0] = 'Cabo Verde'; $arr[112,1] = 5063; $arr[112,2] = 85; $arr[112,3] = 4465; $arr[112,4] = 549; $arr[112,5] = 0; $arr[112,6] = 2; $arr[112,7] = 49
$arr[113,0] = 'Guinea Ecuatorial'; $arr[113,1] = 5002; $arr[113,2] = 2; $arr[113,3] = 4509; $arr[113,4] = 410; $arr[113,5] = 0; $arr[113,6] = 0; $arr[113,7] = 83
$arr[114,0] = 'Hong Kong'; $arr[114,1] = 4994; $arr[114,2] = 9; $arr[114,3] = 4682; $arr[114,4] = 210; $arr[114,5] = 0; $arr[114,6] = 0; $arr[114,7] = 102
$arr[115,0] = 'Nicaragua'; $arr[115,1] = 4961; $arr[115,2] = 0; $arr[115,3] = 2913; $arr[115,4] = 1901; $arr[115,5] = 0; $arr[115,6] = 0; $arr[115,7] = 147
$arr[116,0] = 'Congo'; $arr[116,1] = 4934; $arr[116,2] = 0; $arr[116,3] = 3887; $arr[116,4] = 959; $arr[116,5] = 0; $arr[116,6] = 0; $arr[116,7] = 88
$arr[117,0] = 'Cuba'; $arr[117,1] = 4933; $arr[117,2] = 57; $arr[117,3] = 4230; $arr[117,4] = 594; $arr[117,5] = 0; $arr[117,6] = 0; $arr[117,7] = 109
$arr[118,0] = 'Republica de Africa Central'; $arr[118,1] = 4782; $arr[118,2] = 0; $arr[118,3] = 1830; $arr[118,4] = 2890; $arr[118,5] = 0; $arr[118,6] = 0; $arr[118,7] = 62
$arr[119,0] = 'Ruanda'; $arr[119,1] = 4653; $arr[119,2] = 19; $arr[119,3] = 2817; $arr[119,4] = 1813; $arr[119,5] = 0; $arr[119,6] = 1; $arr[119,7] = 23
$arr[120,0] = 'Surinam'; $arr[120,1] = 4645; $arr[120,2] = 0; $arr[120,3] = 4089; $arr[120,4] = 461; $arr[120,5] = 0; $arr[120,6] = 0; $arr[120,7] = 95
$arr[121,0] = 'Jamaica'; $arr[121,1] = 4374; $arr[121,2] = 210; $arr[121,3] = 1225; $arr[121,4] = 3098; $arr[121,5] = 0; $arr[121,6] = 5; $arr[121,7] = 51
$arr[122,0] = 'Jordania'; $arr[122,1] = 4131; $arr[122,2] = 279; $arr[122,3] = 2415; $arr[122,4] = 1690; $arr[122,5] = 0; $arr[122,6] = 0; $arr[122,7] = 26
$arr[123,0] = 'Eslovenia'; $arr[123,1] = 4058; $arr[123,2] = 104; $arr[123,3] = 2897; $arr[123,4] = 1025; $arr[123,5] = 0; $arr[123,6] = 1; $arr[123,7] = 136
$arr[124,0] = 'Birmania'; $arr[124,1] = 4043; $arr[124,2] = 222; $arr[124,3] = 944; $arr[124,4] = 3039; $arr[124,5] = 0; $arr[124,6] = 20; $arr[124,7] = 60
$arr[125,0] = 'Angola'; $arr[125,1] = 3789; $arr[125,2] = 114; $arr[125,3] = 1405; $arr[125,4] = 2240; $arr[125,5] = 0; $arr[125,6] = 1; $arr[125,7] = 144
$arr[126,0] = 'Siria'; $arr[126,1] = 3691; $arr[126,2] = 37; $arr[126,3] = 903; $arr[126,4] = 2623; $arr[126,5] = 0; $arr[126,6] = 2; $arr[126,7] = 165
$arr[127,0] = 'Lituania'; $arr[127,1] = 3504; $arr[127,2] = 62; $arr[127,3] = 2149; $arr[127,4] = 1268; $arr[127,5] = 0; $arr[127,6] = 0; $arr[127,7] = 87
$arr[128,0] = 'Tailandia'; $arr[128,1] = 3490; $arr[128,2] = 0; $arr[128,3] = 3325; $arr[128,4] = 107; $arr[128,5] = 0; $arr[128,6] = 0; $arr[128,7] = 58
$arr[129,0] = 'Gambia'; $arr[129,1] = 3440; $arr[129,2] = 0; $arr[129,3] = 1851; $arr[129,4] = 1482; $arr[129,5] = 0; $arr[129,6] = 0; $arr[129,7] = 107
$arr[130,0] = 'Trinidad yTobago'; $arr[130,1] = 3434; $arr[130,2] = 107; $arr[130,3] = 1469; $arr[130,4] = 1905; $arr[130,5] = 0; $arr[130,6] = 2; $arr[130,7] = 60
$arr[131,0] = 'Guadalupe'; $arr[131,1] = 3426; $arr[131,2] = 0; $arr[131,3] = 837; $arr[131,4] = 2563; $arr[131,5] = 0; $arr[131,6] = 0; $arr[131,7] = 26
$arr[132,0] = 'Somalia'; $arr[132,1] = 3390; $arr[132,2] = 0; $arr[132,3] = 2812; $arr[132,4] = 480; $arr[132,5] = 0; $arr[132,6] = 0; $arr[132,7] = 98
$arr[133,0] = 'Mayotte'; $arr[133,1] = 3374; $arr[133,2] = 0; $arr[133,3] = 2964; $arr[133,4] = 370; $arr[133,5] = 0; $arr[133,6] = 0; $arr[133,7] = 40
$arr[134,0] = 'Aruba'; $arr[134,1] = 3328; $arr[134,2] = 0; $arr[134,3] = 1676; $arr[134,4] = 1630; $arr[134,5] = 0; $arr[134,6] = 0; $arr[134,7] = 22
$arr[135,0] = 'Sri Lanka'; $arr[135,1] = 3274; $arr[135,2] = 3; $arr[135,3] = 3043; $arr[135,4] = 218; $arr[135,5] = 0; $arr[135,6] = 0; $arr[135,7] = 13
$arr[136,0] = 'Reunion'; $arr[136,1] = 3099; $arr[136,2] = 97; $arr[136,3] = 1794; $arr[136,4] = 1290; $arr[136,5] = 0; $arr[136,6] = 0; $arr[136,7] = 15
$arr[137,0] = 'Bahamas'; $arr[137,1] = 3087; $arr[137,2] = 0; $arr[137,3] = 1533; $arr[137,4] = 1485; $arr[137,5] = 0; $arr[137,6] = 0; $arr[137,7] = 69
$arr[138,0] = 'Mali'; $arr[138,1] = 2966; $arr[138,2] = 0; $arr[138,3] = 2311; $arr[138,4] = 527; $arr[138,5] = 0; $arr[138,6] = 0; $arr[138,7] = 128
$arr[139,0] = 'Georgia'; $arr[139,1] = 2937; $arr[139,2] = 179; $arr[139,3] = 1422; $arr[139,4] = 1496; $arr[139,5] = 0; $arr[139,6] = 0; $arr[139,7] = 19
$arr[140,0] = 'Estonia'; $arr[140,1] = 2778; $arr[140,2] = 22; $arr[140,3] = 2337; $arr[140,4] = 377; $arr[140,5] = 0; $arr[140,6] = 0; $arr[140,7] = 64
$arr[141,0] = 'Malta'; $arr[141,1] = 2595; $arr[141,2] = 35; $arr[141,3] = 1978; $arr[141,4] = 601; $arr[141,5] = 0; $arr[141,6] = 0; $arr[141,7] = 16
$arr[142,0] = 'Sudan del Sur'; $arr[142,1] = 2594; $arr[142,2] = 0; $arr[142,3] = 1290; $arr[142,4] = 1255; $arr[142,5] = 0; $arr[142,6] = 0; $arr[142,7] = 49
$arr[143,0] = 'Botsuana'; $arr[143,1] = 2567; $arr[143,2] = 104; $arr[143,3] = 624; $arr[143,4] = 1930; $arr[143,5] = 0; $arr[143,6] = 2; $arr[143,7] = 13
$arr[144,0] = 'Benin'; $arr[144,1] = 2280; $arr[144,2] = 0; $arr[144,3] = 1942; $arr[144,4] = 298; $arr[144,5] = 0; $arr[144,6] = 0; $arr[144,7] = 40
$arr[145,0] = 'Guinea-Bisau'; $arr[145,1] = 2275; $arr[145,2] = 0; $arr[145,3] = 1127; $arr[145,4] = 1109; $arr[145,5] = 0; $arr[145,6] = 0; $arr[145,7] = 39
$arr[146,0] = 'Islandia'; $arr[146,1] = 2189; $arr[146,2] = 0; $arr[146,3] = 2104; $arr[146,4] = 75; $arr[146,5] = 0; $arr[146,6] = 0; $arr[146,7] = 10
$arr[147,0] = 'Sierra Leona'; $arr[147,1] = 2133; $arr[147,2] = 0; $arr[147,3] = 1646; $arr[147,4] = 415; $arr[147,5] = 0; $arr[147,6] = 0; $arr[147,7] = 72
$arr[148,0] = 'Yemen'; $arr[148,1] = 2022; $arr[148,2] = 3; $arr[148,3] = 1221; $arr[148,4] = 216; $arr[148,5] = 0; $arr[148,6] = 2; $arr[148,7] = 585
$arr[149,0] = 'Guyana'; $arr[149,1] = 1958; $arr[149,2] = 0; $arr[149,3] = 1302; $arr[149,4] = 598; $arr[149,5] = 0; $arr[149,6] = 0; $arr[149,7] = 58
$arr[150,0] = 'Uruguay'; $arr[150,1] = 1856; $arr[150,2] = 0; $arr[150,3] = 1559; $arr[150,4] = 252; $arr[150,5] = 0; $arr[150,6] = 0; $arr[150,7] = 45
$arr[151,0] = 'Nueva Zelanda'; $arr[151,1] = 1809; $arr[151,2] = 7; $arr[151,3] = 1707; $arr[151,4] = 77; $arr[151,5] = 0; $arr[151,6] = 0; $arr[151,7] = 25
$arr[152,0] = 'Burkina Faso'; $arr[152,1] = 1767; $arr[152,2] = 19; $arr[152,3] = 1166; $arr[152,4] = 545; $arr[152,5] = 0; $arr[152,6] = 0; $arr[152,7] = 56
$arr[153,0] = 'Togo'; $arr[153,1] = 1618; $arr[153,2] = 10; $arr[153,3] = 1243; $arr[153,4] = 334; $arr[153,5] = 0; $arr[153,6] = 1; $arr[153,7] = 41
$arr[154,0] = 'Republica de Chipre'; $arr[154,1] = 1558; $arr[154,2] = 10; $arr[154,3] = 1282; $arr[154,4] = 254; $arr[154,5] = 0; $arr[154,6] = 0; $arr[154,7] = 22
$arr[155,0] = 'Belice'; $arr[155,1] = 1536; $arr[155,2] = 8; $arr[155,3] = 696; $arr[155,4] = 821; $arr[155,5] = 0; $arr[155,6] = 0; $arr[155,7] = 19
$arr[156,0] = 'Letonia'; $arr[156,1] = 1494; $arr[156,2] = 8; $arr[156,3] = 1248; $arr[156,4] = 210; $arr[156,5] = 0; $arr[156,6] = 1; $arr[156,7] = 36
$arr[157,0] = 'Principado de Andorra'; $arr[157,1] = 1483; $arr[157,2] = 0; $arr[157,3] = 1054; $arr[157,4] = 376; $arr[157,5] = 0; $arr[157,6] = 0; $arr[157,7] = 53
$arr[158,0] = 'Liberia'; $arr[158,1] = 1333; $arr[158,2] = 1; $arr[158,3] = 1214; $arr[158,4] = 37; $arr[158,5] = 0; $arr[158,6] = 0; $arr[158,7] = 82
$arr[159,0] = 'Lesoto'; $arr[159,1] = 1327; $arr[159,2] = 0; $arr[159,3] = 687; $arr[159,4] = 607; $arr[159,5] = 0; $arr[159,6] = 0; $arr[159,7] = 33
$arr[160,0] = 'Niger'; $arr[160,1] = 1182; $arr[160,2] = 0; $arr[160,3] = 1104; $arr[160,4] = 9; $arr[160,5] = 0; $arr[160,6] = 0; $arr[160,7] = 69
$arr[161,0] = 'Martinica'; $arr[161,1] = 1122; $arr[161,2] = 0; $arr[161,3] = 98; $arr[161,4] = 1006; $arr[161,5] = 0; $arr[161,6] = 0; $arr[161,7] = 18
$arr[162,0] = 'Republica del Chad'; $arr[162,1] = 1115; $arr[162,2] = 25; $arr[162,3] = 962; $arr[162,4] = 72; $arr[162,5] = 0; $arr[162,6] = 0; $arr[162,7] = 81
$arr[163,0] = 'Polinesia Francesa'; $arr[163,1] = 1099; $arr[163,2] = 0; $arr[163,3] = 672; $arr[163,4] = 425; $arr[163,5] = 0; $arr[163,6] = 0; $arr[163,7] = 2
$arr[164,0] = 'Vietnam'; $arr[164,1] = 1066; $arr[164,2] = 3; $arr[164,3] = 940; $arr[164,4] = 91; $arr[164,5] = 0; $arr[164,6] = 0; $arr[164,7] = 35
$arr[165,0] = 'Santo Tome y Principe'; $arr[165,1] = 907; $arr[165,2] = 0; $arr[165,3] = 871; $arr[165,4] = 21; $arr[165,5] = 0; $arr[165,6] = 0; $arr[165,7] = 15
$arr[166,0] = 'San Marino'; $arr[166,1] = 723; $arr[166,2] = 0; $arr[166,3] = 663; $arr[166,4] = 18; $arr[166,5] = 0; $arr[166,6] = 0; $arr[166,7] = 42
$arr[167,0] = 'Crucero'; $arr[167,1] = 712; $arr[167,2] = 0; $arr[167,3] = 651; $arr[167,4] = 48; $arr[167,5] = 0; $arr[167,6] = 0; $arr[167,7] = 13
$arr[168,0] = 'Islas Turcas y Caicos'; $arr[168,1] = 659; $arr[168,2] = 9; $arr[168,3] = 561; $arr[168,4] = 93; $arr[168,5] = 0; $arr[168,6] = 0; $arr[168,7] = 5
$arr[169,0] = 'San Martin (Parte Holandesa)'; $arr[169,1] = 557; $arr[169,2] = 8; $arr[169,3] = 477; $arr[169,4] = 61; $arr[169,5] = 0; $arr[169,6] = 0; $arr[169,7] = 19
$arr[170,0] = 'Papua Nueva Guinea'; $arr[170,1] = 516; $arr[170,2] = 0; $arr[170,3] = 232; $arr[170,4] = 278; $arr[170,5] = 0; $arr[170,6] = 0; $arr[170,7] = 6
$arr[171,0] = 'Tanzania'; $arr[171,1] = 509; $arr[171,2] = 0; $arr[171,3] = 183; $arr[171,4] = 305; $arr[171,5] = 0; $arr[171,6] = 0; $arr[171,7] = 21
$arr[172,0] = 'Taiwan'; $arr[172,1] = 503; $arr[172,2] = 3; $arr[172,3] = 478; $arr[172,4] = 18; $arr[172,5] = 0; $arr[172,6] = 0; $arr[172,7] = 7
$arr[173,0] = 'Burundi'; $arr[173,1] = 472; $arr[173,2] = 0; $arr[173,3] = 374; $arr[173,4] = 97; $arr[173,5] = 0; $arr[173,6] = 0; $arr[173,7] = 1
$arr[174,0] = 'Comoras'; $arr[174,1] = 467; $arr[174,2] = 0; $arr[174,3] = 446; $arr[174,4] = 14; $arr[174,5] = 0; $arr[174,6] = 0; $arr[174,7] = 7
$arr[175,0] = 'Islas Feroe'; $arr[175,1] = 429; $arr[175,2] = 1; $arr[175,3] = 412; $arr[175,4] = 17; $arr[175,5] = 0; $arr[175,6] = 0; $arr[175,7] = 0
$arr[176,0] = 'Mauricio'; $arr[176,1] = 366; $arr[176,2] = 1; $arr[176,3] = 338; $arr[176,4] = 18; $arr[176,5] = 0; $arr[176,6] = 0; $arr[176,7] = 10
$arr[177,0] = 'Eritrea'; $arr[177,1] = 364; $arr[177,2] = 0; $arr[177,3] = 305; $arr[177,4] = 59; $arr[177,5] = 0; $arr[177,6] = 0; $arr[177,7] = 0
$arr[178,0] = 'Gibraltar'; $arr[178,1] = 343; $arr[178,2] = 3; $arr[178,3] = 313; $arr[178,4] = 30; $arr[178,5] = 0; $arr[178,6] = 0; $arr[178,7] = 0
$arr[179,0] = 'Isla de Man'; $arr[179,1] = 339; $arr[179,2] = 0; $arr[179,3] = 312; $arr[179,4] = 3; $arr[179,5] = 0; $arr[179,6] = 0; $arr[179,7] = 24
$arr[180,0] = 'San Martin (Parte Francesa)'; $arr[180,1] = 330; $arr[180,2] = 0; $arr[180,3] = 206; $arr[180,4] = 118; $arr[180,5] = 0; $arr[180,6] = 0; $arr[180,7] = 6
$arr[181,0] = 'Mongolia'; $arr[181,1] = 311; $arr[181,2] = 0; $arr[181,3] = 301; $arr[181,4] = 10; $arr[181,5] = 0; $arr[181,6] = 0; $arr[181,7] = 0
$arr[182,0] = 'Camboya'; $arr[182,1] = 275; $arr[182,2] = 0; $arr[182,3] = 274; $arr[182,4] = 1; $arr[182,5] = 0; $arr[182,6] = 0; $arr[182,7] = 0
$arr[183,0] = 'Butan'; $arr[183,1] = 246; $arr[183,2] = 0; $arr[183,3] = 175; $arr[183,4] = 71; $arr[183,5] = 0; $arr[183,6] = 0; $arr[183,7] = 0
$arr[184,0] = 'Curazao'; $arr[184,1] = 210; $arr[184,2] = 18; $arr[184,3] = 74; $arr[184,4] = 135; $arr[184,5] = 0; $arr[184,6] = 0; $arr[184,7] = 1
$arr[185,0] = 'Islas Caimanes'; $arr[185,1] = 208; $arr[185,2] = 0; $arr[185,3] = 204; $arr[185,4] = 3; $arr[185,5] = 0; $arr[185,6] = 0; $arr[185,7] = 1
$arr[186,0] = 'Monaco'; $arr[186,1] = 186; $arr[186,2] = 5; $arr[186,3] = 147; $arr[186,4] = 38; $arr[186,5] = 0; $arr[186,6] = 0; $arr[186,7] = 1
$arr[187,0] = 'Barbados'; $arr[187,1] = 185; $arr[187,2] = 0; $arr[187,3] = 171; $arr[187,4] = 7; $arr[187,5] = 0; $arr[187,6] = 0; $arr[187,7] = 7
$arr[188,0] = 'Bermudas'; $arr[188,1] = 178; $arr[188,2] = 0; $arr[188,3] = 164; $arr[188,4] = 5; $arr[188,5] = 0; $arr[188,6] = 0; $arr[188,7] = 9
$arr[189,0] = 'Brunei'; $arr[189,1] = 145; $arr[189,2] = 0; $arr[189,3] = 141; $arr[189,4] = 1; $arr[189,5] = 0; $arr[189,6] = 0; $arr[189,7] = 3
$arr[190,0] = 'Seychelles'; $arr[190,1] = 141; $arr[190,2] = 1; $arr[190,3] = 136; $arr[190,4] = 5; $arr[190,5] = 0; $arr[190,6] = 0; $arr[190,7] = 0
$arr[191,0] = 'Liechtenstein'; $arr[191,1] = 112; $arr[191,2] = 1; $arr[191,3] = 108; $arr[191,4] = 3; $arr[191,5] = 0; $arr[191,6] = 0; $arr[191,7] = 1
$arr[192,0] = 'Antigua y Barbuda'; $arr[192,1] = 95; $arr[192,2] = 0; $arr[192,3] = 91; $arr[192,4] = 1; $arr[192,5] = 0; $arr[192,6] = 0; $arr[192,7] = 3
$arr[193,0] = 'Islas Virgenes Britanicas'; $arr[193,1] = 66; $arr[193,2] = 0; $arr[193,3] = 37; $arr[193,4] = 28; $arr[193,5] = 0; $arr[193,6] = 0; $arr[193,7] = 1
$arr[194,0] = 'San Vicente y las Granadinas'; $arr[194,1] = 64; $arr[194,2] = 0; $arr[194,3] = 64; $arr[194,4] = 0; $arr[194,5] = 0; $arr[194,6] = 0; $arr[194,7] = 0
$arr[195,0] = 'Macao'; $arr[195,1] = 46; $arr[195,2] = 0; $arr[195,3] = 46; $arr[195,4] = 0; $arr[195,5] = 0; $arr[195,6] = 0; $arr[195,7] = 0
$arr[196,0] = 'Puerto Rico'; $arr[196,1] = 39; $arr[196,2] = 0; $arr[196,3] = 1; $arr[196,4] = 36; $arr[196,5] = 0; $arr[196,6] = 0; $arr[196,7] = 2
$arr[197,0] = 'Guam'; $arr[197,1] = 32; $arr[197,2] = 0; $arr[197,3] = 0; $arr[197,4] = 31; $arr[197,5] = 0; $arr[197,6] = 0; $arr[197,7] = 1
$arr[198,0] = 'Bonaire, San Eustaquio y Saba'; $arr[198,1] = 32; $arr[198,2] = 4; $arr[198,3] = 17; $arr[198,4] = 14; $arr[198,5] = 0; $arr[198,6] = 0; $arr[198,7] = 1
$arr[199,0] = 'Fiyi'; $arr[199,1] = 32; $arr[199,2] = 0; $arr[199,3] = 24; $arr[199,4] = 6; $arr[199,5] = 0; $arr[199,6] = 0; $arr[199,7] = 2
$arr[200,0] = 'Santa Lucia'; $arr[200,1] = 27; $arr[200,2] = 0; $arr[200,3] = 26; $arr[200,4] = 1; $arr[200,5] = 0; $arr[200,6] = 0; $arr[200,7] = 0
$arr[201,0] = 'Timor Oriental'; $arr[201,1] = 27; $arr[201,2] = 0; $arr[201,3] = 26; $arr[201,4] = 1; $arr[201,5] = 0; $arr[201,6] = 0; $arr[201,7] = 0
$arr[202,0] = 'Nueva Caledonia'; $arr[202,1] = 26; $arr[202,2] = 0; $arr[202,3] = 26; $arr[202,4] = 0; $arr[202,5] = 0; $arr[202,6] = 0; $arr[202,7] = 0
$arr[203,0] = 'Dominica'; $arr[203,1] = 24; $arr[203,2] = 0; $arr[203,3] = 18; $arr[203,4] = 6; $arr[203,5] = 0; $arr[203,6] = 0; $arr[203,7] = 0
$arr[204,0] = 'Granada'; $arr[204,1] = 24; $arr[204,2] = 0; $arr[204,3] = 24; $arr[204,4] = 0; $arr[204,5] = 0; $arr[204,6] = 0; $arr[204,7] = 0
$arr[205,0] = 'San Bartolome'; $arr[205,1] = 23; $arr[205,2] = 0; $arr[205,3] = 16; $arr[205,4] = 7; $arr[205,5] = 0; $arr[205,6] = 0; $arr[205,7] = 0
$arr[206,0] = 'Laos'; $arr[206,1] = 23; $arr[206,2] = 0; $arr[206,3] = 22; $arr[206,4] = 1; $arr[206,5] = 0; $arr[206,6] = 0; $arr[206,7] = 0
$arr[207,0] = 'Islas Virgenes de los Estados Unidos'; $arr[207,1] = 17; $arr[207,2] = 0; $arr[207,3] = 0; $arr[207,4] = 17; $arr[207,5] = 0; $arr[207,6] = 0; $arr[207,7] = 0
$arr[208,0] = 'San Cristobal y Nieves'; $arr[208,1] = 17; $arr[208,2] = 0; $arr[208,3] = 17; $arr[208,4] = 0; $arr[208,5] = 0; $arr[208,6] = 0; $arr[208,7] = 0
$arr[209,0] = 'Groenlandia'; $arr[209,1] = 14; $arr[209,2] = 0; $arr[209,3] = 14; $arr[209,4] = 0; $arr[209,5] = 0; $arr[209,6] = 0; $arr[209,7] = 0
$arr[210,0] = 'Montserrat'; $arr[210,1] = 13; $arr[210,2] = 0; $arr[210,3] = 12; $arr[210,4] = 0; $arr[210,5] = 0; $arr[210,6] = 0; $arr[210,7] = 1
$arr[211,0] = 'Islas Malvinas'; $arr[211,1] = 13; $arr[211,2] = 0; $arr[211,3] = 13; $arr[211,4] = 0; $arr[211,5] = 0; $arr[211,6] = 0; $arr[211,7] = 0
$arr[212,0] = 'Santa Sede'; $arr[212,1] = 12; $arr[212,2] = 0; $arr[212,3] = 12; $arr[212,4] = 0; $arr[212,5] = 0; $arr[212,6] = 0; $arr[212,7] = 0
$arr[213,0] = 'San Pedro y Miquelon'; $arr[213,1] = 11; $arr[213,2] = 0; $arr[213,3] = 5; $arr[213,4] = 6; $arr[213,5] = 0; $arr[213,6] = 0; $arr[213,7] = 0
$arr[214,0] = 'Sahara Occidental'; $arr[214,1] = 10; $arr[214,2] = 0; $arr[214,3] = 8; $arr[214,4] = 1; $arr[214,5] = 0; $arr[214,6] = 0; $arr[214,7] = 1
$arr[215,0] = 'Anguila'; $arr[215,1] = 3; $arr[215,2] = 0; $arr[215,3] = 3; $arr[215,4] = 0; $arr[215,5] = 0; $arr[215,6] = 0; $arr[215,7] = 0

$ws.Range("A4:H219").Value = $arr

Write-Output "Updated country table: $rows rows (A4:H219), refreshed header date."
